# Update the "as of" disclaimer date from 2021-05-18 to 2021-05-19,
# and refresh the Weight (D) / Percent Change (E) figures for the
# holdings rows (2-35) to the newly-reported values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Unprotect()

# --- Disclaimer text: bump the "as of" date by one day -------------------
$oldText = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution." + [char]10 + "Model holdings provided as of 2021-05-18 for illustrative purposes only and are subject to change."
$newText = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution." + [char]10 + "Model holdings provided as of 2021-05-19 for illustrative purposes only and are subject to change."

$found = $ws.Cells.Find($oldText)
if ($found) {
    $found.Value = $newText
} else {
    $ws.Range("A38").Value = $newText
}

# --- Refreshed Weight (D) / Percent Change (E) values ---------------------
$ws.Range("D2").Value = 0.03558057102378977
$ws.Range("E2").Value = -0.000394788787998257
$ws.Range("D3").Value = 0.02021576105261792
$ws.Range("E3").Value = 0.003153330705557567
$ws.Range("D4").Value = 0.01937908193928529
$ws.Range("E4").Value = -0.0008012820512820484
$ws.Range("D5").Value = 0.03763727947912179
$ws.Range("E5").Value = -0.0003528581510233719
$ws.Range("D6").Value = 0.03435083661271023
$ws.Range("E6").Value = -0.000400000000000067
$ws.Range("D7").Value = 0.01977183076116652
$ws.Range("E7").Value = -0.00253435327489604
$ws.Range("D8").Value = 0.03721919531901021
$ws.Range("E8").Value = -0.006807062327164326
$ws.Range("D9").Value = 0.02044071433803409
$ws.Range("E9").Value = -0.004677941705649502
$ws.Range("D10").Value = 0.02617559289543994
$ws.Range("E10").Value = -0.0129768757927603
$ws.Range("D11").Value = 0.02407751019823983
$ws.Range("E11").Value = -0.007690267833465936
$ws.Range("D12").Value = 0.05740308288159642
$ws.Range("E12").Value = -0.005464480874317057
$ws.Range("D13").Value = 0.02473260236115136
$ws.Range("E13").Value = -0.0007434944237918462
$ws.Range("D14").Value = 0.02687343844161517
$ws.Range("E14").Value = -0.002805486284289338
$ws.Range("D15").Value = 0.032712008000246
$ws.Range("E15").Value = -0.004547044421126301
$ws.Range("D16").Value = 0.01984354611373789
$ws.Range("E16").Value = -0.005493133583021281
$ws.Range("D17").Value = 0.03129373985245638
$ws.Range("E17").Value = -0.00612094997143553
$ws.Range("D18").Value = 0.0419161933572582
$ws.Range("E18").Value = -0.001389210465385426
$ws.Range("D19").Value = 0.125724470637007
$ws.Range("E19").Value = 0
$ws.Range("D20").Value = 0.00932575411706908
$ws.Range("E20").Value = -0.0002519526329048949
$ws.Range("D21").Value = 0.01531873643269652
$ws.Range("E21").Value = -0.001470485260135934
$ws.Range("D22").Value = 0.01728013089379906
$ws.Range("E22").Value = -0.01170854271356792
$ws.Range("D23").Value = 0.01545021457907736
$ws.Range("E23").Value = -0.004684684684684748
$ws.Range("D24").Value = 0.02142598316691523
$ws.Range("E24").Value = -0.002767811378780016
$ws.Range("D25").Value = 0.01272758514691278
$ws.Range("E25").Value = -0.02151115891368649
$ws.Range("D26").Value = 0.0423925079318508
$ws.Range("E26").Value = -0.004446131323537328
$ws.Range("D27").Value = 0.0239617134003201
$ws.Range("E27").Value = 0
$ws.Range("D28").Value = 0.04558675264200073
$ws.Range("E28").Value = -0.004773269689737347
$ws.Range("D29").Value = 0.05527373964609597
$ws.Range("E29").Value = -0.0003622532149972768
$ws.Range("D30").Value = 0.01313422754137194
$ws.Range("E30").Value = -0.01355713363460298
$ws.Range("D31").Value = 0.02064830065772785
$ws.Range("E31").Value = 0.0003834355828220559
$ws.Range("D32").Value = 0.01328700576041821
$ws.Range("E32").Value = 0.0004805382027870397
$ws.Range("D33").Value = 0.04183027795624321
$ws.Range("E33").Value = -0.001549586776859568
$ws.Range("D34").Value = 0.01700961486301703
$ws.Range("E34").Value = -0.01088875809299583
$ws.Range("D35").Value = 0.9999999999999999
$ws.Range("E35").Value = -0.003304320625243551

# --- Restore the sheet protection that was in force before the edit ------
$ws.Protect("", $true, $true, $true, $false, $false, $true, $true)
